$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("2025")
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = 46436.48618400006
$ws1.Range("E2").Value = 254562.6397561083
$ws1.Range("G2").Value = 64767.40570129472
$ws1.Range("I2").Value = 169331.8557977695
$ws1.Range("L2").Value = 391439.272052782
$ws1.Range("N2").Value = 57497.45683913826
$ws1.Range("O2").Value = 52873.29882886782

$ws2 = $wb.Worksheets.Item("2030")
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = 53213.57539827293
$ws2.Range("E2").Value = 140913.2167893911
$ws2.Range("I2").Value = 137718.8074294309
$ws2.Range("L2").Value = 100616.5925057371
$ws2.Range("N2").Value = 13951.98301190023
$ws2.Range("O2").Value = 5784.902224332255

$ws3 = $wb.Worksheets.Item("2035")
$ws3.Range("A2").Value = 9466.406440477318
$ws3.Range("B2").Value = 28163.43894958384
$ws3.Range("E2").Value = 134040.3978993126
$ws3.Range("I2").Value = 142850.152056096
$ws3.Range("L2").Value = 0
$ws3.Range("M2").Value = 0
$ws3.Range("N2").Value = 35153.79306694151
$ws3.Range("O2").Value = 39339.82394484724
